$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Expected Wins")

$values = @(
    -0.166666666666667,
    -0.666666666666667,
    1.083333333333333,
    -0.5,
    -0.5833333333333335,
    -0.75,
    0.1666666666666665,
    0.9166666666666665,
    -0.5,
    -0.3333333333333333,
    0.4166666666666667,
    0.9166666666666666
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Range("D$row").Value = $values[$i]
}
